# Add 2022-Q3 data:
#  - insert a new sheet "2022-Q3" right after "总计", seeded as a copy of the
#    previous-quarter sheet (so header styling / column layout matches),
#    then overwrite its single data row with the new quarter's fund info
#    and drop the rest of the old rows.
#  - record the new quarter in the "总计" roll-up sheet (new row 2, pushing
#    every existing row down by one).

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item(1)
$prevQuarter = $wb.Worksheets.Item(2)

# --- 1. New "2022-Q3" sheet, placed right after "总计" ------------------
$prevQuarter.Copy($null, $total)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# Drop the old quarter's extra fund rows (2..11 had data; keep header+row2)
$q3.Range("3:11").Delete()

# Overwrite row 2 with the single 2022-Q3 fund holding
$q3.Range("B2:G2").NumberFormat = "@"
$q3.Cells.Item(2,2).Value = "513690"
$q3.Cells.Item(2,3).Value = "博时恒生港股通高股息率ETF"
$q3.Cells.Item(2,4).Value = "3.05"
$q3.Cells.Item(2,5).Value = "97.26"
$q3.Cells.Item(2,6).Value = "2.64"
$q3.Cells.Item(2,7).Value = "0.0805"
$q3.Cells.Item(2,8).Value = 6

# --- 2. Add the 2022-Q3 summary row to "总计" ---------------------------
$total.Rows.Item(2).Insert()
$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q3"
$total.Cells.Item(2,3).Value = 1
$total.Cells.Item(2,4).Value = 0.08
$total.Range("B2:D2").ClearFormats()
$total.Cells.Item(3,1).Copy()
$total.Cells.Item(2,1).PasteSpecial(-4122)

# Re-number the index column (A) for the rows that got pushed down
$total.Cells.Item(3,1).Value = 1
$total.Cells.Item(4,1).Value = 2
$total.Cells.Item(5,1).Value = 3
$total.Cells.Item(6,1).Value = 4
$total.Cells.Item(7,1).Value = 5
$total.Cells.Item(8,1).Value = 6
$total.Cells.Item(9,1).Value = 7
